# Auto-generated: apply scheduled market-price refresh to profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 712.05884
$ws.Range("I6").Value = 21
$ws.Range("K6").Value = 63
$ws.Range("M6").Value = 49
$ws.Range("H15").Value = 530.4375
$ws.Range("I15").Value = 530.4375
$ws.Range("K15").Value = 1591.3125
$ws.Range("M15").Value = -1422.3125
$ws.Range("H28").Value = 697.1724
$ws.Range("I28").Value = 768.92
$ws.Range("J28").Value = 248.75
$ws.Range("K28").Value = 768.92
$ws.Range("L28").Value = 248.75
$ws.Range("M28").Value = -283.92
$ws.Range("N28").Value = -1218.75
$ws.Range("H32").Value = 4606.727
$ws.Range("I32").Value = 4579.5
$ws.Range("K32").Value = 4579.5
$ws.Range("M32").Value = -4253.5
$ws.Range("H43").Value = 11498251
$ws.Range("J43").Value = 4602.227
$ws.Range("L43").Value = 4602.227
$ws.Range("N43").Value = -4740.227
$ws.Range("H132").Value = 22354.666
$ws.Range("I132").Value = 25518.158
$ws.Range("J132").Value = 2469.8572
$ws.Range("K132").Value = 76554.474
$ws.Range("L132").Value = 7409.571599999999
$ws.Range("M132").Value = -74024.474
$ws.Range("N132").Value = -12469.5716
$ws.Range("H137").Value = 1924672.4
$ws.Range("I137").Value = 967.37036
$ws.Range("J137").Value = 4002273.8
$ws.Range("K137").Value = 2902.11108
$ws.Range("L137").Value = 12006821.4
$ws.Range("M137").Value = -352.1110800000001
$ws.Range("N137").Value = -12011921.4
$ws.Range("H139").Value = 107994.14
$ws.Range("J139").Value = 124659.836
$ws.Range("L139").Value = 124659.836
$ws.Range("N139").Value = -134939.836

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12271219
$ws.Range("I32").Value = 14935106
$ws.Range("K32").Value = 14935106
$ws.Range("M32").Value = -14934819
$ws.Range("H33").Value = 100000000
$ws.Range("I33").Value = 100000000
$ws.Range("K33").Value = 100000000
$ws.Range("M33").Value = -99999671
$ws.Range("H36").Value = 6266.6665
$ws.Range("I36").Value = 9000
$ws.Range("J36").Value = 800
$ws.Range("K36").Value = 9000
$ws.Range("L36").Value = 800
$ws.Range("M36").Value = -8654
$ws.Range("N36").Value = -1492
$ws.Range("H61").Value = 2186006.8
$ws.Range("I61").Value = 7829
$ws.Range("K61").Value = 7829
$ws.Range("M61").Value = -7617
$ws.Range("H132").Value = 2327.625
$ws.Range("I132").Value = 2382.0588
$ws.Range("J132").Value = 2195.4285
$ws.Range("K132").Value = 7146.176399999999
$ws.Range("L132").Value = 6586.2855
$ws.Range("M132").Value = -4616.176399999999
$ws.Range("N132").Value = -11646.2855
$ws.Range("H136").Value = 2186006.8
$ws.Range("I136").Value = 7829
$ws.Range("K136").Value = 23487
$ws.Range("M136").Value = -20937

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2921.04
$ws.Range("J86").Value = 4888
$ws.Range("L86").Value = 4888
$ws.Range("N86").Value = -7134
$ws.Range("H89").Value = 2921.04
$ws.Range("J89").Value = 4888
$ws.Range("L89").Value = 24440
$ws.Range("N89").Value = -35672
$ws.Range("H134").Value = 17310536
$ws.Range("I134").Value = 2486.1428
$ws.Range("J134").Value = 52944756
$ws.Range("K134").Value = 7458.428400000001
$ws.Range("L134").Value = 158834268
$ws.Range("M134").Value = -4923.428400000001
$ws.Range("N134").Value = -158839338
$ws.Range("H138").Value = 146660.83
$ws.Range("J138").Value = 146660.83
$ws.Range("L138").Value = 146660.83
$ws.Range("N138").Value = -156940.83

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1016.6667
$ws.Range("I2").Value = 1140
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 1140
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -1027
$ws.Range("N2").Value = -626
$ws.Range("H3").Value = 1764922.5
$ws.Range("I3").Value = 3509845
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 3509845
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = -3509732
$ws.Range("N3").Value = -20226
$ws.Range("H22").Value = 601.26666
$ws.Range("I22").Value = 636.6667
$ws.Range("K22").Value = 636.6667
$ws.Range("M22").Value = -286.6667
$ws.Range("H31").Value = 4078.2114
$ws.Range("I31").Value = 1870.5834
$ws.Range("J31").Value = 4740.5
$ws.Range("K31").Value = 1870.5834
$ws.Range("L31").Value = 4740.5
$ws.Range("M31").Value = -1575.5834
$ws.Range("N31").Value = -5330.5
$ws.Range("H34").Value = 4078.2114
$ws.Range("I34").Value = 1870.5834
$ws.Range("J34").Value = 4740.5
$ws.Range("K34").Value = 1870.5834
$ws.Range("L34").Value = 4740.5
$ws.Range("M34").Value = -1668.5834
$ws.Range("N34").Value = -5144.5
$ws.Range("H134").Value = 2114.1853
$ws.Range("I134").Value = 1603.2609
$ws.Range("K134").Value = 4809.7827
$ws.Range("M134").Value = -2274.7827

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 968.1818
$ws.Range("J5").Value = 1600
$ws.Range("L5").Value = 4800
$ws.Range("N5").Value = -5024
$ws.Range("H23").Value = 125037.625
$ws.Range("J23").Value = 166703.5
$ws.Range("L23").Value = 500110.5
$ws.Range("N23").Value = -500580.5
$ws.Range("H70").Value = 12119
$ws.Range("J70").Value = 12472.167
$ws.Range("L70").Value = 37416.501
$ws.Range("N70").Value = -38046.501
$ws.Range("H73").Value = 12119
$ws.Range("J73").Value = 12472.167
$ws.Range("L73").Value = 37416.501
$ws.Range("N73").Value = -39600.501
$ws.Range("H95").Value = 21663.334
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 21663.334
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 64990.00199999999
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -69108.00199999999
$ws.Range("H107").Value = 83334040
$ws.Range("J107").Value = 200000260
$ws.Range("L107").Value = 600000780
$ws.Range("N107").Value = -600004620
$ws.Range("H113").Value = 1010.4545
$ws.Range("J113").Value = 1000.58826
$ws.Range("L113").Value = 3001.76478
$ws.Range("N113").Value = -7341.76478
$ws.Range("H132").Value = 2055.7666
$ws.Range("J132").Value = 2138.4
$ws.Range("L132").Value = 19245.6
$ws.Range("N132").Value = -24305.6
$ws.Range("H135").Value = 968.1818
$ws.Range("J135").Value = 1600
$ws.Range("L135").Value = 14400
$ws.Range("N135").Value = -19470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 11999
$ws.Range("J24").Value = 11999
$ws.Range("L24").Value = 11999
$ws.Range("N24").Value = -12345
$ws.Range("H97").Value = 4452.4614
$ws.Range("I97").Value = 670.875
$ws.Range("K97").Value = 670.875
$ws.Range("M97").Value = -174.875
$ws.Range("H113").Value = 2616.3547
$ws.Range("I113").Value = 1935.1428
$ws.Range("J113").Value = 4046.9
$ws.Range("K113").Value = 1935.1428
$ws.Range("L113").Value = 4046.9
$ws.Range("M113").Value = 234.8571999999999
$ws.Range("N113").Value = -8386.9
$ws.Range("H126").Value = 4230
$ws.Range("J126").Value = 4955.5
$ws.Range("L126").Value = 14866.5
$ws.Range("N126").Value = -19806.5
$ws.Range("H132").Value = 9848476
$ws.Range("I132").Value = 4044.6667
$ws.Range("J132").Value = 21207434
$ws.Range("K132").Value = 12134.0001
$ws.Range("L132").Value = 63622302
$ws.Range("M132").Value = -9604.000100000001
$ws.Range("N132").Value = -63627362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 550000000
$ws.Range("J20").Value = 550000000
$ws.Range("L20").Value = 550000000
$ws.Range("N20").Value = -550000452
$ws.Range("H43").Value = 3039999.8
$ws.Range("J43").Value = 3039999.8
$ws.Range("L43").Value = 3039999.8
$ws.Range("N43").Value = -3040385.8
$ws.Range("H46").Value = 8786.571
$ws.Range("I46").Value = 26105.75
$ws.Range("K46").Value = 26105.75
$ws.Range("M46").Value = -25917.75
$ws.Range("H68").Value = 5888.75
$ws.Range("I68").Value = 7185.3335
$ws.Range("K68").Value = 7185.3335
$ws.Range("M68").Value = -6436.3335
$ws.Range("H71").Value = 5888.75
$ws.Range("I71").Value = 7185.3335
$ws.Range("K71").Value = 35926.6675
$ws.Range("M71").Value = -32182.6675
$ws.Range("H129").Value = 68500
$ws.Range("J129").Value = 68500
$ws.Range("L129").Value = 68500
$ws.Range("N129").Value = -78500
$ws.Range("H136").Value = 4512
$ws.Range("I136").Value = 3799.923
$ws.Range("J136").Value = 5129.1333
$ws.Range("K136").Value = 11399.769
$ws.Range("L136").Value = 15387.3999
$ws.Range("M136").Value = -8849.769
$ws.Range("N136").Value = -20487.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 39999.75
$ws.Range("J64").Value = 39999.75
$ws.Range("L64").Value = 39999.75
$ws.Range("N64").Value = -40495.75
$ws.Range("H67").Value = 39999.75
$ws.Range("J67").Value = 39999.75
$ws.Range("L67").Value = 39999.75
$ws.Range("N67").Value = -41715.75
$ws.Range("H107").Value = 58824070
$ws.Range("J107").Value = 125000510
$ws.Range("L107").Value = 375001530
$ws.Range("N107").Value = -375005370
$ws.Range("H126").Value = 2436.35
$ws.Range("I126").Value = 2265.1538
$ws.Range("K126").Value = 6795.4614
$ws.Range("M126").Value = -4325.4614
$ws.Range("H129").Value = 49950
$ws.Range("J129").Value = 49950
$ws.Range("L129").Value = 49950
$ws.Range("N129").Value = -59950
$ws.Range("H132").Value = 72961.42999999999
$ws.Range("I132").Value = 112275.555
$ws.Range("K132").Value = 336826.665
$ws.Range("M132").Value = -334296.665
